# Updates the crypto price/volume table on Sheet1 to the latest scraped
# values (GitHub Actions refresh). Most rows only get their Price (D) and
# Volume(1h) (E) columns refreshed; two pairs of rows (30/31 and 45/46)
# also swap coin name/link because their relative ranking changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that *looks* numeric (e.g. "24.70", "1.00") while
# forcing the cell to stay plain text, matching the source data which
# stores every Price/Volume cell as a string (not a number).
function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
}

# --- Row 2 (Bitcoin) ---
$ws.Range('D2').Value = '60.424.35'
$ws.Range('E2').Value = '  -0.72%  '

# --- Row 3 (Ethereum) ---
$ws.Range('D3').Value = '2.612.56'
$ws.Range('E3').Value = '  +0.21%  '

# --- Row 4 (TetherUSD) ---
$ws.Range('E4').Value = '  +0.02%  '

# --- Row 5 (BNB) ---
Set-TextCell 'D5' '582.58'
$ws.Range('E5').Value = '  +2.07%  '

# --- Row 6 (Solana) ---
Set-TextCell 'D6' '143.28'
$ws.Range('E6').Value = '  -0.09%  '

# --- Row 7 (USDC) ---
Set-TextCell 'D7' '0.998'
$ws.Range('E7').Value = '  +0.27%  '

# --- Row 8 (XRP) ---
$ws.Range('E8').Value = '  -0.87%  '

# --- Row 9 (Toncoin) ---
Set-TextCell 'D9' '6.51'
$ws.Range('E9').Value = '  -0.21%  '

# --- Row 10 (Dogecoin) ---
$ws.Range('E10').Value = '  -1.51%  '

# --- Row 11 (TRON) ---
Set-TextCell 'D11' '0.156'
$ws.Range('E11').Value = '  +1.10%  '

# --- Row 12 (Cardano) ---
Set-TextCell 'D12' '0.374'
$ws.Range('E12').Value = '  +1.27%  '

# --- Row 13 (WrappedliquidstakedEther2.0) ---
$ws.Range('D13').Value = '3.069.65'
$ws.Range('E13').Value = '  -0.51%  '

# --- Row 14 (Avalanche) ---
Set-TextCell 'D14' '24.70'
$ws.Range('E14').Value = '  +4.96%  '

# --- Row 15 (WrappedBTC) ---
$ws.Range('D15').Value = '60.400.17'
$ws.Range('E15').Value = '  -0.71%  '

# --- Row 16 (ShibaInu) ---
$ws.Range('E16').Value = '  -0.07%  '

# --- Row 17 (WrappedEther) ---
$ws.Range('D17').Value = '2.613.31'
$ws.Range('E17').Value = '  -0.27%  '

# --- Row 18 (Chainlink) ---
Set-TextCell 'D18' '11.37'
$ws.Range('E18').Value = '  +0.62%  '

# --- Row 19 (Polkadot) ---
$ws.Range('E19').Value = '  -0.87%  '

# --- Row 20 (BitcoinCash) ---
Set-TextCell 'D20' '346.57'
$ws.Range('E20').Value = '  -0.35%  '

# --- Row 21 (Uniswap) ---
Set-TextCell 'D21' '6.93'
$ws.Range('E21').Value = '  -2.61%  '

# --- Row 22 (Dai) ---
Set-TextCell 'D22' '0.998'
$ws.Range('E22').Value = '  -0.60%  '

# --- Row 23 (Polygon) ---
Set-TextCell 'D23' '0.533'
$ws.Range('E23').Value = '  +2.01%  '

# --- Row 24 (Litecoin) ---
Set-TextCell 'D24' '63.59'
$ws.Range('E24').Value = '  -0.70%  '

# --- Row 25 (Binance-PegBSC-USD) ---
Set-TextCell 'D25' '1.00'
$ws.Range('E25').Value = '  +0.52%  '

# --- Row 26 (Kaspa) ---
$ws.Range('E26').Value = '  -0.01%  '

# --- Row 27 (InternetComputer(DFINITY)) ---
$ws.Range('E27').Value = '  +3.22%  '

# --- Row 28 (PancakeSwap) ---
$ws.Range('E28').Value = '  +5.13%  '

# --- Row 29 (PEPE) ---
$ws.Range('E29').Value = '  -0.08%  '

# --- Row 30/31: Monero and Aptos swap ranking order ---
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 'D30' '6.44'
$ws.Range('E30').Value = '  +1.82%  '

$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell 'D31' '168.67'
$ws.Range('E31').Value = '  +4.59%  '

# --- Row 32 (USDe) ---
$ws.Range('E32').Value = '  +0.20%  '

# --- Row 33 (EthereumClassic) ---
$ws.Range('E33').Value = '  -0.28%  '

# --- Row 34 (ImmutableX) ---
$ws.Range('E34').Value = '  +8.97%  '

# --- Row 35 (NEARProtocol) ---
$ws.Range('E35').Value = '  +0.31%  '

# --- Row 36 (Fetch.AI) ---
$ws.Range('E36').Value = '  +4.53%  '

# --- Row 37 (Stacks) ---
$ws.Range('E37').Value = '  +3.89%  '

# --- Row 38 (Bittensor) ---
Set-TextCell 'D38' '319.67'

# --- Row 39 (OKB) ---
Set-TextCell 'D39' '38.33'
$ws.Range('E39').Value = '  +1.57%  '

# --- Row 40 (Filecoin) ---
$ws.Range('E40').Value = '  +2.56%  '

# --- Row 41 (SuiNetwork) ---
Set-TextCell 'D41' '0.852'
$ws.Range('E41').Value = '  -0.71%  '

# --- Row 42 (Aave) ---
Set-TextCell 'D42' '135.59'
$ws.Range('E42').Value = '  -2.80%  '

# --- Row 43 (Stellar) ---
Set-TextCell 'D43' '0.0993'
$ws.Range('E43').Value = '  +0.51%  '

# --- Row 44 (FirstDigitalUSD) ---
$ws.Range('E44').Value = '  +0.23%  '

# --- Row 45/46: EnergySwap and Mantle swap ranking order ---
$ws.Range('B45').Value = 'Mantle'
$ws.Range('C45').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell 'D45' '0.611'
$ws.Range('E45').Value = '  +0.54%  '

$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D46' '19.87'
$ws.Range('E46').Value = '  +1.23%  '

# --- Row 47 (RenderToken) ---
Set-TextCell 'D47' '5.03'
$ws.Range('E47').Value = '  +4.40%  '

# --- Row 48 (Hedera) ---
$ws.Range('E48').Value = '  -0.46%  '

# --- Row 49 (InjectiveProtocol) ---
$ws.Range('E49').Value = '  +1.56%  '

# --- Row 50 (VeChain) ---
$ws.Range('E50').Value = '  -0.33%  '

# --- Row 51 (WhiteBITCoin) ---
Set-TextCell 'D51' '10.74'
$ws.Range('E51').Value = '  +0.29%  '
